# issue #5: stock data from json to db
#
# Adds three metadata columns to the "股票" (stock) worksheet:
#   - "category"    inserted right after the existing "property_category"
#                    column (holds "normal")
#   - "source_file" appended after "legislator_id" (holds "tmp4c4f1")
#   - "index"       appended after "source_file" (holds the row's original
#                    disclosure index, same value as column A)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# ---- Header row (row 1) ----------------------------------------------
# Existing headers B1:K1 are untouched; I1 slides from "date" to the new
# "category" header, and J1/K1 become "date"/"legislator_name". L1:N1 are
# brand-new trailing headers.
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"

$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:N1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# ---- Data rows 2-4 ------------------------------------------------------
$indexValues = @{ 2 = 79; 3 = 80; 4 = 81 }

foreach ($r in @(2, 3, 4)) {
    $ws.Range("H$r").Value = "stock"
    $ws.Range("I$r").Value = "normal"

    # "2012-04-26" must stay a literal text string, not get auto-parsed
    # into a date serial number; force text format right before the write.
    $ws.Range("J$r").NumberFormat = "@"
    $ws.Range("J$r").Value = "2012-04-26"

    $ws.Range("K$r").Value = "黃昭順"
    $ws.Range("L$r").Value = 665
    $ws.Range("M$r").Value = "tmp4c4f1"
    $ws.Range("N$r").Value = $indexValues[$r]

    # normalise formatting on the touched/new cells back to the plain
    # data-row look (this also drops the temporary "@" text format on J,
    # without disturbing any of the values just written).
    $ws.Range("C$r").Copy() | Out-Null
    $ws.Range("H$r`:N$r").PasteSpecial(-4122) | Out-Null
    $ws.Application.CutCopyMode = $false
}
